# "Levels are displayed hierarchically"
#
# The request-parameter table (row 16, under the "LVL / field name / data
# type / max length / required / field name & description" header in row
# 15) and the response-parameter table (row 23, under the header in row
# 22) used to hold a single literal "level" placeholder in column A.
#
# Going forward the level column holds the *value* for a specific row
# (e.g. "1", "1.1", "1.1.1", ...) so the placeholder is renamed from
# "level" to "levelValue", and the column is switched to a plain Text
# number format (so Excel never "helpfully" reinterprets a value like
# "1.10" as a number and drops the trailing zero) while keeping the same
# border/font it already had.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Request parameter table: A16 ---------------------------------------
$reqLevel = $ws.Range("A16")
$reqLevel.Value2 = '${requestParameter.levelValue}'
$reqLevel.NumberFormat = "@"

# --- Response parameter table: A23 --------------------------------------
$respLevel = $ws.Range("A23")
$respLevel.Value2 = '${responseParameter.levelValue}'
$respLevel.NumberFormat = "@"

# Reflect where the author was last working when the sheet was saved.
$ws.Range("A23").Select()
$excel.ActiveWindow.ScrollRow = 13
